# Weekly fruit/vegetable price update: insert two new records (Primera /
# Segunda quality rows) for the latest reporting week at the top of the
# existing Cilantro price history block, pushing the previously-first rows
# (and everything below them) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before row 287 (shifts rows 287:307 -> 289:309).
$ws.Range("A287:A288").EntireRow.Insert()

# --- New row 287: "Primera" quality ---
$ws.Cells.Item(287, 1).Value = 7
$ws.Cells.Item(287, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(287, 3).Value = "Ñuble"
$ws.Cells.Item(287, 4).Value = 45166
$ws.Cells.Item(287, 5).Value = 16
$ws.Cells.Item(287, 6).Value = 100112040
$ws.Cells.Item(287, 7).Value = "Cilantro"
$ws.Cells.Item(287, 8).Value = "Sin especificar"
$ws.Cells.Item(287, 9).Value = "Primera"
$ws.Cells.Item(287, 10).Value = 250
$ws.Cells.Item(287, 11).Value = 1500
$ws.Cells.Item(287, 12).Value = 1500
$ws.Cells.Item(287, 13).Value = 1500
$ws.Cells.Item(287, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(287, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(287, 16).Value = 1500
$ws.Cells.Item(287, 17).Value = 1
$ws.Cells.Item(287, 18).Value = "Hortaliza"

# --- New row 288: "Segunda" quality ---
$ws.Cells.Item(288, 1).Value = 7
$ws.Cells.Item(288, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(288, 3).Value = "Ñuble"
$ws.Cells.Item(288, 4).Value = 45166
$ws.Cells.Item(288, 5).Value = 16
$ws.Cells.Item(288, 6).Value = 100112040
$ws.Cells.Item(288, 7).Value = "Cilantro"
$ws.Cells.Item(288, 8).Value = "Sin especificar"
$ws.Cells.Item(288, 9).Value = "Segunda"
$ws.Cells.Item(288, 10).Value = 200
$ws.Cells.Item(288, 11).Value = 1000
$ws.Cells.Item(288, 12).Value = 1000
$ws.Cells.Item(288, 13).Value = 1000
$ws.Cells.Item(288, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(288, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(288, 16).Value = 1000
$ws.Cells.Item(288, 17).Value = 1
$ws.Cells.Item(288, 18).Value = "Hortaliza"
